# daily auto push: 2026-02-24 14:15 UTC
# Insert two new rows of data (date 2026/02/24) right before the existing
# 2026/12/29 block, which starts at row 850. Inserting shifts all the rows
# that used to be 850..891 down to 852..893, which is exactly what the
# target workbook looks like (dimension grows from D891 to D893).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 850.. down by two to make room for the new entries.
$ws.Rows("850:851").Insert()

function Set-TextCell($range, [string]$text) {
    # Writing a date-shaped string straight to .Value lets Excel's COM
    # layer "helpfully" reinterpret it as a date serial. Force the cell to
    # be read as plain text first, then restore the default (no explicit
    # number-format) styling so the final cell matches a normal, unstyled
    # inline-string cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# New row 850: 2026/02/24, 火, 20:00, rank 46
Set-TextCell $ws.Range("A850") "2026/02/24"
$ws.Range("B850").Value = "火"
$ws.Range("C850").Value = 20
$ws.Range("D850").Value = 46

# New row 851: 2026/02/24, 火, 22:00, rank 45
Set-TextCell $ws.Range("A851") "2026/02/24"
$ws.Range("B851").Value = "火"
$ws.Range("C851").Value = 22
$ws.Range("D851").Value = 45
